$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 - InceptionResNetV2_10
$ws.Range("A42").Value = "InceptionResNetV2_10"
$ws.Range("B42").Value = "InceptionResNetV2"
$ws.Range("C42").Value = 0.9678391959798995
$ws.Range("D42").Value = 0.92
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = "['Functional', 'Dropout', 'BatchNormalization', 'Dense', 'Dropout', 'BatchNormalization', 'Dense']"
$ws.Range("G42").Value = "[0.1, 0.1]"
$ws.Range("H42").Value = "[('relu', 1024), ('softmax', 15)]"
$ws.Range("I42").Value = 70
$ws.Range("J42").Value = 0.0001
$ws.Range("K42").Value = "{'monitor': 'val_loss', 'patience': 15, 'min_delta': 0, 'restore_best_weights': True}"
$ws.Range("L42").Value = "{'Train': 32, 'Validation': 32}"
$ws.Range("M42").Value = "{'zoom_range': 0.2, 'rotation_range': 30, 'shear_range': 0.2, 'brightness_range': None, 'horizontal_flip': True, 'width_shift_range': 0.2, 'height_shift_range': 0.2}"
$ws.Range("N42").Value = 55936239
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 2777.928627490997

# Row 43 - InceptionResNetV2_18
$ws.Range("A43").Value = "InceptionResNetV2_18"
$ws.Range("B43").Value = "InceptionResNetV2"
$ws.Range("C43").Value = 0.9447236180904522
$ws.Range("D43").Value = 0.9166666666666666
$ws.Range("E43").Value = 7
$ws.Range("F43").Value = "['Functional', 'Dropout', 'BatchNormalization', 'Dense', 'Dropout', 'BatchNormalization', 'Dense']"
$ws.Range("G43").Value = "[0.1, 0.1]"
$ws.Range("H43").Value = "[('relu', 1024), ('softmax', 15)]"
$ws.Range("I43").Value = 70
$ws.Range("J43").Value = 0.00001
$ws.Range("K43").Value = "{'monitor': 'val_loss', 'patience': 15, 'min_delta': 0, 'restore_best_weights': True}"
$ws.Range("L43").Value = "{'Train': 32, 'Validation': 32}"
$ws.Range("M43").Value = "{'zoom_range': 0.2, 'rotation_range': 30, 'shear_range': 0.2, 'brightness_range': None, 'horizontal_flip': True, 'width_shift_range': 0.2, 'height_shift_range': 0.2}"
$ws.Range("N43").Value = 55936239
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 3645.254849910736
